$d = $word.ActiveDocument

$replacements = @(
    @("711÷6=", "923÷8="),
    @("174÷7=", "925÷2="),
    @("755÷3=", "797÷4="),
    @("647÷6=", "407÷5="),
    @("333÷8=", "422÷7="),
    @("599÷5=", "826÷7="),
    @("847÷3=", "813÷5="),
    @("978÷2=", "159÷8="),
    @("840÷5=", "700÷6="),
    @("563÷8=", "809÷9="),
    @("977÷3=", "825÷6="),
    @("766÷8=", "245÷9="),
    @("351÷5=", "686÷6="),
    @("292÷3=", "109÷2="),
    @("754÷6=", "642÷2="),
    @("758÷4=", "114÷4="),
    @("976÷9=", "499÷9="),
    @("629÷7=", "566÷3="),
    @("410÷9=", "899÷9="),
    @("530÷7=", "937÷6="),
    @("496÷4=", "670÷7="),
    @("442÷5=", "275÷2="),
    @("779÷8=", "808÷9="),
    @("840÷9=", "803÷6="),
    @("489÷8=", "726÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
